$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from an existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J (rows 2-10)
$values = @(
    @(6, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(6, 7),
    @(6, 7),
    @(7, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
